$d = $word.ActiveDocument

# 1) Merge the two runs that spell "...ci-d" + bookmark + "essus." into a
#    single run reading "...ci-dessus." (the stray mid-word bookmark goes
#    away here and reappears below, next to the caisse manager's score).
$d.Content.Find.Execute(
    "décrits ci-dessus.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "décrits ci-dessus.",
    2
) | Out-Null

# 2) Locate the evaluation table and recolor/re-bookmark the relevant cells.
$t = $d.Tables.Item(1)

# Row 6 = "Paiement : les étapes décrites sont implantées (WC)" -> "Fonction" cell (holds "4")
$monnayeurCell = $t.Cell(6, 2)

# Put the _GoBack bookmark back, collapsed, right before the "4".
$bmRange = $monnayeurCell.Range
$bmRange.Collapse(1) # wdCollapseStart
$bmRange = $d.Range($bmRange.Start, $bmRange.Start)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# Recolor that cell from yellow to green (caisse manager service is now done).
$monnayeurCell.Shading.Texture = 0
$monnayeurCell.Shading.ForegroundPatternColor = -16777216
$monnayeurCell.Shading.BackgroundPatternColor = 5287936

# Row 7 = "Paiement : une option pour voir l'état de la caisse est au menu (WC)"
# -> "Fonction" cell currently has no shading at all; give it the same green.
$caisseCell = $t.Cell(7, 2)
$caisseCell.Shading.Texture = 0
$caisseCell.Shading.ForegroundPatternColor = -16777216
$caisseCell.Shading.BackgroundPatternColor = 5287936
